# Sanity.xlsx - "TestCaseMaster" sheet
# 1. changed the folder structure (new FileName path for the Multiple
#    Allotment cases: Sanity//Accural//MultipleLeaveAllotment.xlsx)
# 2. added cases for the missed "Multiple Allotment" scenarios
#    (rows 164 & 165: TCID 169 "...WithEncashment" and
#     TCID 170 "...WithLeaves")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 164: Multiple Allotment With Encashment - Financial ---
$ws.Range("A164").Value = "169"
$ws.Range("D164").Value = "com.darwinbox.leaves.Accural.MultipleAllotment.Daily.Fiancial.MutlipleAllotmentWithEncashment"
$ws.Range("E164").Value = "Sanity//Accural//MultipleLeaveAllotment.xlsx"
$ws.Range("B164").Value = "Multiple Allotment With Encashment- Financial"
$ws.Range("C164").Value = "Multiple Allotment With Encashment- Financial"
$ws.Range("F164").Value = "Financial"
$ws.Range("G164").Value = "All"

# --- Row 165: Multiple Allotment With Leaves - Financial ---
$ws.Range("A165").Value = "170"
$ws.Range("B165").Value = "Multiple Allotment With Leaves- Financial"
$ws.Range("C165").Value = "Multiple Allotment With Leaves- Financial"
$ws.Range("D165").Value = "com.darwinbox.leaves.Accural.MultipleAllotment.Daily.Fiancial.MutlipleAllotmentWithLeaves"
$ws.Range("E165").Value = "Sanity//Accural//MultipleLeaveAllotment.xlsx"
$ws.Range("F165").Value = "Financial"
$ws.Range("G165").Value = "All"

# Match the author's final selection/scroll position (last row touched)
[void]$ws.Range("A165").Select()
